# Apply updated cryptocurrency price/volume data from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new price strings look like plain numbers (e.g. "569.26") and Excel would
# otherwise auto-convert them to a numeric cell. Force those specific cells to Text
# format first so they are stored as strings (matching the source data), then restore
# each cell default styling once its text value is committed.
$textForceCells = @('D5', 'D6', 'D10', 'D11', 'D12', 'D14', 'D19', 'D21', 'D22', 'D23', 'D24', 'D26', 'D29', 'D30', 'D31', 'D32', 'D33', 'D35', 'D36', 'D37', 'D39', 'D41', 'D42', 'D43', 'D45', 'D46', 'D48', 'D50', 'D51')
foreach ($addr in $textForceCells) { $ws.Range($addr).NumberFormat = "@" }

# -- Price column (D) updates (number-like strings) --
$ws.Range('D5').Value = '569.26'
$ws.Range('D6').Value = '141.47'
$ws.Range('D10').Value = '7.50'
$ws.Range('D11').Value = '0.124'
$ws.Range('D12').Value = '0.394'
$ws.Range('D14').Value = '28.47'
$ws.Range('D19').Value = '6.21'
$ws.Range('D21').Value = '8.97'
$ws.Range('D22').Value = '383.23'
$ws.Range('D23').Value = '0.559'
$ws.Range('D24').Value = '73.89'
$ws.Range('D26').Value = '0.0000116'
$ws.Range('D29').Value = '0.999'
$ws.Range('D30').Value = '7.41'
$ws.Range('D31').Value = '7.99'
$ws.Range('D32').Value = '2.14'
$ws.Range('D33').Value = '1.43'
$ws.Range('D35').Value = '23.59'
$ws.Range('D36').Value = '6.98'
$ws.Range('D37').Value = '165.32'
$ws.Range('D39').Value = '4.98'
$ws.Range('D41').Value = '28.18'
$ws.Range('D42').Value = '0.0772'
$ws.Range('D43').Value = '1.00'
$ws.Range('D45').Value = '41.96'
$ws.Range('D46').Value = '4.42'
$ws.Range('D48').Value = '1.13'
$ws.Range('D50').Value = '23.46'
$ws.Range('D51').Value = '6.81'

# Restore original (default) formatting now that the text values are committed.
foreach ($addr in $textForceCells) { $ws.Range($addr).Style = "Normal" }

# -- Remaining cells (already safe as text, no numeric auto-conversion risk) --
$ws.Range('D2').Value = '60.821.54'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '3.388.47'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -1.62%  '
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.388.56'
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '3.964.64'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').Value = '3.385.37'
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').Value = '60.905.25'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('E21').Value = '  -5.93%  '
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  -5.59%  '
$ws.Range('D27').Value = '3.522.28'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -2.72%  '
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E33').Value = '  -2.52%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('D38').Value = '3.417.18'
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('E39').Value = '  -2.78%  '
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('E47').Value = '  -4.06%  '
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('D49').Value = '2.488.24'
$ws.Range('E49').Value = '  -4.35%  '
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('E51').Value = '  -1.50%  '
